# Generate Report for Handoff
# Updates the localization-status report: the "32b7cbc2..." file was handed
# back and is replaced in the report by a newly-generated handoff entry
# "d5e47a05-7745-416b-a80f-2a61e8b68b7a.md", and the "e763c5cb..." file
# (now "ffffa777b6d6-87e7-4a44-a1ec-7a147d912d17.md") moves from
# "Handed back: in sync with en-US" to "Ready for handoff" / "True" to
# be localized, with new handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldFile1 = "32b7cbc2-45b2-451f-8e8d-1ab2cc0211a5.md"
$newFile1 = "d5e47a05-7745-416b-a80f-2a61e8b68b7a.md"
$oldFile2 = "e763c5cb-5d0e-42a3-a729-ac235ed9e9fd.md"
$newFile2 = "ffffa777b6d6-87e7-4a44-a1ec-7a147d912d17.md"

$newXlfZhCn = "d5e47a05-7745-416b-a80f-2a61e8b68b7a.0fb762ba73f259f3995538a124d3941c824ea16e.zh-cn.xlf"
$newXlfDeDe = "d5e47a05-7745-416b-a80f-2a61e8b68b7a.0fb762ba73f259f3995538a124d3941c824ea16e.de-de.xlf"

$statusText = "Ready for handoff"
$handoffDate = "2016-08-18 23:04:20"
$zhCnGenDate = "2016-08-18 23:04:14"
$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newFile1
$ws1.Range("A3").Value = $newFile2

$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText

$ws1.Range("G2").Value = $handoffDate
$ws1.Range("G3").Value = $handoffDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile1", "", "", "e2e\$newFile1")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile2", "", "", "e2e\$newFile2")

$ws1.Columns.Item(5).ColumnWidth = 16.25
$ws1.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newFile1
$ws2.Range("A3").Value = $newFile2

$ws2.Range("C2").Value = $statusText
$ws2.Range("C3").Value = $statusText

$ws2.Range("F3").Value = "True"

$ws2.Range("G2").Value = $newXlfZhCn
$ws2.Range("G3").Value = $newXlfZhCn

$ws2.Range("H2").Value = $zhCnGenDate
$ws2.Range("H3").Value = $zhCnGenDate

$ws2.Range("I2").Value = ""
$ws2.Range("I3").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("I3").Style = "Normal"

$ws2.Range("J2").Value = ""
$ws2.Range("J3").Value = ""

$ws2.Range("K2").Value = $zeroDate
$ws2.Range("K3").Value = $zeroDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile1", "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile2", "", "", $newFile2)

$ws2.Columns.Item(3).ColumnWidth = 16.25
$ws2.Columns.Item(9).ColumnWidth = 17.75
$ws2.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newFile1
$ws3.Range("A3").Value = $newFile2

$ws3.Range("C2").Value = $statusText
$ws3.Range("C3").Value = $statusText

$ws3.Range("F3").Value = "True"

$ws3.Range("G2").Value = $newXlfDeDe
$ws3.Range("G3").Value = $newXlfDeDe

$ws3.Range("H2").Value = $handoffDate
$ws3.Range("H3").Value = $handoffDate

$ws3.Range("I2").Value = ""
$ws3.Range("I3").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("I3").Style = "Normal"

$ws3.Range("J2").Value = ""
$ws3.Range("J3").Value = ""

$ws3.Range("K2").Value = $zeroDate
$ws3.Range("K3").Value = $zeroDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile1", "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea7d55a6eb73e723419466ffb2cc5f4b451d623/e2e/$newFile2", "", "", $newFile2)

$ws3.Columns.Item(3).ColumnWidth = 16.25
$ws3.Columns.Item(9).ColumnWidth = 17.75
$ws3.Columns.Item(10).ColumnWidth = 20.75
